$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) "...is an implementation procedurally generating a city using..."
#    -> "...is an implementation that can procedurally generate a city using..."
# -----------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute(
    "is an implementation procedurally generating a city using",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "is an implementation that can procedurally generate a city using",
    2)

# -----------------------------------------------------------------------
# 2) "...the cities generated a user study will be conducted."
#    -> "...the cities generated, a user study will be conducted."
# -----------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute(
    "the cities generated a user study will be conducted",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the cities generated, a user study will be conducted",
    2)

# -----------------------------------------------------------------------
# 3) "...recognized: Districts, Roads along with blocks and individual
#    houses. All three stages..." ->
#    "...recognized: Districts, roads (creating blocks) and buildings.
#    All three stages..." with "Districts", "roads" and "buildings"
#    italicised.
# -----------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute(
    "three different generation stages have been recognized: Districts, Roads along with blocks and individual houses. All three stages will be procedurally generated with",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "three different generation stages have been recognized: Districts, roads (creating blocks) and buildings. All three stages will be procedurally generated with",
    2)

$scope = $d.Content
$scope.Find.Execute(
    "recognized: Districts, roads (creating blocks) and buildings. All",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rDistricts = $d.Range($scope.Start, $scope.End)
$rDistricts.Find.Execute("Districts", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rDistricts.Italic = 1

$rRoads = $d.Range($rDistricts.End, $scope.End)
$rRoads.Find.Execute("roads", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rRoads.Italic = 1

$rBuildings = $d.Range($rRoads.End, $scope.End)
$rBuildings.Find.Execute("buildings", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rBuildings.Italic = 1

# -----------------------------------------------------------------------
# 4) "...noise in the following order: Districts, roads with blocks and
#    lastly houses. Four different..." ->
#    "...noise in the following order: Districts, roads and lastly
#    building. Four different..."
# -----------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute(
    "noise in the following order: Districts, roads with blocks and lastly houses. Four different",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "noise in the following order: Districts, roads and lastly building. Four different",
    2)

# -----------------------------------------------------------------------
# 5) Move the "_GoBack" bookmark from the end of the "4.4 Viability"
#    heading paragraph to the end of the "...parameters are entered. "
#    paragraph (right after the trailing space, before the paragraph
#    mark).
# -----------------------------------------------------------------------
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

$r5 = $d.Content
$r5.Find.Execute(
    "if all the parameters are entered. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $r5.End

# Use a tiny round-trip marker so we can obtain a non-zero-length Range
# (zero-length Range/Bookmarks.Add at this exact paragraph-end boundary
# is unreliable) precisely at the end of the paragraph text, then shrink
# it back down to a true collapsed bookmark.
$rLast = $d.Range($endPos - 1, $endPos)
$rLast.Text = " ZZMARKZZ"

$r6 = $d.Content
$r6.Find.Execute("ZZMARKZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $r6)

$r7 = $d.Content
$r7.Find.Execute("ZZMARKZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r7.Text = ""
